# Update header labels on row 1 of each sheet so Power BI can treat the
# first row as an automatic header.
$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5 and 6 use the "Ano ####" pattern.
# Sheet 4 ("Potencia Incremental - SIN(MW)") uses the "Intervalo ..." pattern.
$anoSheets = @(1, 2, 3, 5, 6)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value()
    if ($idx -ne 6) {
        $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Value()
        $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Value()
        $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Value()
    }
}

$wsIntervalo = $wb.Worksheets.Item(4)
$wsIntervalo.Range("B1").Value = "Intervalo " + $wsIntervalo.Range("B1").Value()
$wsIntervalo.Range("C1").Value = "Intervalo " + $wsIntervalo.Range("C1").Value()
$wsIntervalo.Range("D1").Value = "Intervalo " + $wsIntervalo.Range("D1").Value()
$wsIntervalo.Range("E1").Value = "Intervalo " + $wsIntervalo.Range("E1").Value()
